# Update the build timestamp embedded in the version string throughout the workbook.
$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

$oldVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $oldTimestamp)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..."
$wsAbout.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation text containing the version string
$a6Text = [string]$wsAbout.Range("A6").Value2
$a6Text = $a6Text.Replace($oldVersion, $newVersion)
$wsAbout.Range("A6").Value = $a6Text

# Column S, rows 2-10 on the data sheet hold the plain version string
for ($row = 2; $row -le 10; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
